# metals_prices.xlsx update: new weekly price pull (12/04/2024) + new "CUSN" function/sheet.
#
# Strategy for literal text: every data cell in this workbook stores plain TEXT
# (dates like "12/04/2024" and decimals like "29,025" are literal strings, not
# real Excel dates/numbers). Assigning such look-alike strings straight to
# .Value lets Excel's input-parser reinterpret them (date serials, or
# thousands-grouped numbers). Forcing NumberFormat="@" before the write avoids
# that reinterpretation, and ClearFormats() afterwards drops the now-unneeded
# number-format override so the cell is left with no explicit style, matching
# every other plain data cell in the sheet.
function Set-TextCell {
    param($Cell, [string]$Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.ClearFormats()
}

function Add-PriceRows {
    param($ws, [int]$startRow, [int]$count, [string]$colA, [string]$colB, [string]$colC, [string]$colD)
    for ($i = 0; $i -lt $count; $i++) {
        $r = $startRow + $i
        Set-TextCell $ws.Cells.Item($r, 1) $colA
        Set-TextCell $ws.Cells.Item($r, 2) $colB
        Set-TextCell $ws.Cells.Item($r, 3) $colC
        Set-TextCell $ws.Cells.Item($r, 4) $colD
    }
}

$wb = $excel.ActiveWorkbook

# --- Per-metal sheets: each gets two more rows of the new (12/04/2024 /
#     "Semaine 15") quote, repeating the same values (matches the source
#     feed's habit of writing the same scraped row twice). ---

$ws = $wb.Worksheets.Item("2360")
Add-PriceRows $ws 29 2 "Semaine 15" "22,59" "€" "KG"

$ws = $wb.Worksheets.Item("1AG1")
Add-PriceRows $ws 99 2 "12/04/2024" "1 051,71 " "€" "KG"

$ws = $wb.Worksheets.Item("1AG2")
Add-PriceRows $ws 97 2 "12/04/2024" "29,025" "$" "OZ"

$ws = $wb.Worksheets.Item("3AL1")
Add-PriceRows $ws 95 2 "12/04/2024" "2443,00" "$" "TO"

$ws = $wb.Worksheets.Item("1AU2")
Add-PriceRows $ws 95 2 "12/04/2024" "2401,5" "$" "OZ"

$ws = $wb.Worksheets.Item("1AU3")
Add-PriceRows $ws 96 2 "12/04/2024" "72 665,00" "€" "KG"

$ws = $wb.Worksheets.Item("2B16")
Add-PriceRows $ws 94 2 "12/04/2024" "1238,00" "€" "100KG"

$ws = $wb.Worksheets.Item("3CU1")
Add-PriceRows $ws 120 2 "12/04/2024" "9402,00" "$" "TO"

$ws = $wb.Worksheets.Item("3CU3")
Add-PriceRows $ws 94 2 "12/04/2024" "1032,41" "€" "100KG"

$ws = $wb.Worksheets.Item("2CUB")
Add-PriceRows $ws 20 2 "Semaine 15" "10,07" "€" "KG"

$ws = $wb.Worksheets.Item("2M30")
Add-PriceRows $ws 129 2 "12/04/2024" "869,00" "€" "100KG"

$ws = $wb.Worksheets.Item("2M37")
Add-PriceRows $ws 94 2 "12/04/2024" "844,00" "€" "100KG"

$ws = $wb.Worksheets.Item("3NI1")
Add-PriceRows $ws 109 2 "12/04/2024" "17780,00" "$" "TO"
$ws.Columns.Item(1).ColumnWidth = 13

$ws = $wb.Worksheets.Item("3SN1")
Add-PriceRows $ws 113 2 "12/04/2024" "32975,00" "$" "TO"

$ws = $wb.Worksheets.Item("3ZN1")
Add-PriceRows $ws 112 2 "12/04/2024" "2848,00" "$" "TO"

$ws = $wb.Worksheets.Item("ZLME")
Add-PriceRows $ws 15 2 "12/04/2024" "1,06510" "USD" "EUR"

$ws = $wb.Worksheets.Item("EURX")
Add-PriceRows $ws 15 2 "12/04/2024" "1,06520" "USD" "EUR"

# --- RPA: the daily roll-up sheet. Refresh the date/weekday header and every
#     metal's latest quote, then append the new CUSN line. ---

$rpa = $wb.Worksheets.Item("RPA")
Set-TextCell $rpa.Cells.Item(1, 1) "12/04/2024"
Set-TextCell $rpa.Cells.Item(1, 2) "Vendredi"

Set-TextCell $rpa.Cells.Item(2, 3) "22,59"
Set-TextCell $rpa.Cells.Item(3, 3) "1 051,71 "
Set-TextCell $rpa.Cells.Item(4, 3) "29,025"
Set-TextCell $rpa.Cells.Item(5, 3) "2443,00"
Set-TextCell $rpa.Cells.Item(6, 3) "2401,5"
Set-TextCell $rpa.Cells.Item(7, 3) "72 665,00"
Set-TextCell $rpa.Cells.Item(8, 3) "1238,00"
Set-TextCell $rpa.Cells.Item(9, 3) "9402,00"
Set-TextCell $rpa.Cells.Item(10, 3) "1032,41"
Set-TextCell $rpa.Cells.Item(11, 3) "10,07"
Set-TextCell $rpa.Cells.Item(12, 3) "869,00"
Set-TextCell $rpa.Cells.Item(13, 3) "844,00"
Set-TextCell $rpa.Cells.Item(14, 3) "17780,00"
Set-TextCell $rpa.Cells.Item(15, 3) "32975,00"
Set-TextCell $rpa.Cells.Item(16, 3) "2848,00"
Set-TextCell $rpa.Cells.Item(17, 3) "1,06510"
Set-TextCell $rpa.Cells.Item(18, 3) "1,06520"

Set-TextCell $rpa.Cells.Item(19, 1) "CUSN"
Set-TextCell $rpa.Cells.Item(19, 2) "CUSN"
Set-TextCell $rpa.Cells.Item(19, 3) "11681,41"
Set-TextCell $rpa.Cells.Item(19, 4) "€"
Set-TextCell $rpa.Cells.Item(19, 5) "TO"

# --- New "CUSN" sheet: the freshly added metal function, placed after RPA
#     (last tab) and left as the active sheet, as in the authored workbook. ---

$cusn = $wb.Worksheets.Add($null, $rpa)
$cusn.Name = "CUSN"

Set-TextCell $cusn.Cells.Item(2, 1) "Semaine 15"
Set-TextCell $cusn.Cells.Item(2, 3) "€"
Set-TextCell $cusn.Cells.Item(2, 4) "KG"

Set-TextCell $cusn.Cells.Item(3, 1) "12/04/2024"
Set-TextCell $cusn.Cells.Item(3, 2) "Valeur non trouvée"
Set-TextCell $cusn.Cells.Item(3, 3) "€"
Set-TextCell $cusn.Cells.Item(3, 4) "KG"

Set-TextCell $cusn.Cells.Item(4, 1) "12/04/2024"
Set-TextCell $cusn.Cells.Item(4, 2) "11681,41"
Set-TextCell $cusn.Cells.Item(4, 3) "€"
Set-TextCell $cusn.Cells.Item(4, 4) "KG"

Set-TextCell $cusn.Cells.Item(5, 1) "12/04/2024"
Set-TextCell $cusn.Cells.Item(5, 2) "11681,41"
Set-TextCell $cusn.Cells.Item(5, 3) "€"
Set-TextCell $cusn.Cells.Item(5, 4) "TO"

$cusn.Activate()
$cusn.Range("R33").Select()
